$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the cell values as described by the diff
$ws.Range("A1").Value = "tushar"
$ws.Range("B5").Value = "Kamthe"

# Reflect the active selection on the sheet (C6) seen in the diff
$ws.Range("C6").Select()
